# The presentation ships two DrawingML theme parts: theme1.xml (the
# "Office Theme" palette, wired to the Notes Master) and theme2.xml (the
# "Integral" palette, wired to the Slide Master / the presentation's main
# look). The edit swaps the two themes' colour schemes, so the deck's
# visible design switches from "Integral" to the default "Office Theme"
# palette.
#
# The only element that actually differs between the two theme parts is
# <a:clrScheme> (the font/format schemes are identical), so re-pointing
# the presentation's live theme colour scheme to the Office Theme RGB
# values reproduces the swap.

$p = $ppt.ActivePresentation

# Office Theme palette (previously theme1.xml), as COM RGB integers
# (0x00BBGGRR), in ThemeColorScheme.Item(index) order: dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink.
$officeThemeColors = @(
    0x000000,   # 1  dk1      000000
    0xFFFFFF,   # 2  lt1      FFFFFF
    0x6A5444,   # 3  dk2      44546A
    0xE6E6E7,   # 4  lt2      E7E6E6
    0xD59B5B,   # 5  accent1  5B9BD5
    0x317DED,   # 6  accent2  ED7D31
    0xA5A5A5,   # 7  accent3  A5A5A5
    0x00C0FF,   # 8  accent4  FFC000
    0xC47244,   # 9  accent5  4472C4
    0x47AD70,   # 10 accent6  70AD47
    0xC16305,   # 11 hlink    0563C1
    0x724F95    # 12 folHlink 954F72
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeThemeColors[$i - 1]
}
